$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 3 for columns D, M, N, O, P, S
$cols = @("D", "M", "N", "O", "P", "S")
foreach ($col in $cols) {
    $addr2 = "${col}2"
    $addr3 = "${col}3"
    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}

# Swap the data between row 4 and row 5 for columns D, M
$cols2 = @("D", "M")
foreach ($col in $cols2) {
    $addr4 = "${col}4"
    $addr5 = "${col}5"
    $val4 = $ws.Range($addr4).Value2
    $val5 = $ws.Range($addr5).Value2
    $ws.Range($addr4).Value2 = $val5
    $ws.Range($addr5).Value2 = $val4
}
